# Apply weekly update to the "Berenjena" price sheet.
# A new data point was inserted at row 43 (shifting the former rows 43-52
# down to 44-52), so update each row's D, I, J, K, L, M, N, O, P, Q columns
# to the values shown by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 43 - new entry
$ws.Range("D43").Value2 = 45001
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value2 = 40
$ws.Range("K43").Value2 = 10000
$ws.Range("L43").Value2 = 10000
$ws.Range("M43").Value2 = 10000
$ws.Range("N43").Value = "$/caja 60 unidades"
$ws.Range("O43").Value = "Región de Arica y Parinacota"
$ws.Range("P43").Value2 = 167
$ws.Range("Q43").Value2 = 60

# Row 44
$ws.Range("D44").Value2 = 44266
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value2 = 60
$ws.Range("K44").Value2 = 9000
$ws.Range("L44").Value2 = 9500
$ws.Range("M44").Value2 = 9208
$ws.Range("N44").Value = "$/caja 60 unidades"
$ws.Range("O44").Value = "Región del Maule"
$ws.Range("P44").Value2 = 153
$ws.Range("Q44").Value2 = 60

# Row 45
$ws.Range("D45").Value2 = 44812
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value2 = 60
$ws.Range("K45").Value2 = 12000
$ws.Range("L45").Value2 = 13000
$ws.Range("M45").Value2 = 12500
$ws.Range("N45").Value = "$/caja 60 unidades"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value2 = 208
$ws.Range("Q45").Value2 = 60

# Row 46
$ws.Range("D46").Value2 = 44812
$ws.Range("I46").Value = "Segunda"
$ws.Range("J46").Value2 = 60
$ws.Range("K46").Value2 = 14000
$ws.Range("L46").Value2 = 14000
$ws.Range("M46").Value2 = 14000
$ws.Range("N46").Value = "$/caja 90 unidades"
$ws.Range("O46").Value = "Región de Arica y Parinacota"
$ws.Range("P46").Value2 = 156
$ws.Range("Q46").Value2 = 90

# Row 47
$ws.Range("D47").Value2 = 44790
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value2 = 60
$ws.Range("K47").Value2 = 12000
$ws.Range("L47").Value2 = 13000
$ws.Range("M47").Value2 = 12500
$ws.Range("N47").Value = "$/caja 60 unidades"
$ws.Range("O47").Value = "Región de Arica y Parinacota"
$ws.Range("P47").Value2 = 208
$ws.Range("Q47").Value2 = 60

# Row 48
$ws.Range("D48").Value2 = 44594
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value2 = 80
$ws.Range("K48").Value2 = 12000
$ws.Range("L48").Value2 = 13000
$ws.Range("M48").Value2 = 12500
$ws.Range("N48").Value = "$/caja 60 unidades"
$ws.Range("O48").Value = "Región de Arica y Parinacota"
$ws.Range("P48").Value2 = 208
$ws.Range("Q48").Value2 = 60

# Row 49
$ws.Range("D49").Value2 = 44819
$ws.Range("I49").Value = "Segunda"
$ws.Range("J49").Value2 = 60
$ws.Range("K49").Value2 = 14000
$ws.Range("L49").Value2 = 14000
$ws.Range("M49").Value2 = 14000
$ws.Range("N49").Value = "$/caja 90 unidades"
$ws.Range("O49").Value = "Región de Arica y Parinacota"
$ws.Range("P49").Value2 = 156
$ws.Range("Q49").Value2 = 90

# Row 50
$ws.Range("D50").Value2 = 44798
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value2 = 40
$ws.Range("K50").Value2 = 12000
$ws.Range("L50").Value2 = 12000
$ws.Range("M50").Value2 = 12000
$ws.Range("N50").Value = "$/caja 60 unidades"
$ws.Range("O50").Value = "Región de Arica y Parinacota"
$ws.Range("P50").Value2 = 200
$ws.Range("Q50").Value2 = 60

# Row 51
$ws.Range("D51").Value2 = 44764
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value2 = 60
$ws.Range("K51").Value2 = 12000
$ws.Range("L51").Value2 = 13000
$ws.Range("M51").Value2 = 12500
$ws.Range("N51").Value = "$/caja 60 unidades"
$ws.Range("O51").Value = "Región de Arica y Parinacota"
$ws.Range("P51").Value2 = 208
$ws.Range("Q51").Value2 = 60

# Row 52
$ws.Range("D52").Value2 = 44782
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value2 = 60
$ws.Range("K52").Value2 = 12000
$ws.Range("L52").Value2 = 13000
$ws.Range("M52").Value2 = 12500
$ws.Range("N52").Value = "$/caja 60 unidades"
$ws.Range("O52").Value = "Región de Arica y Parinacota"
$ws.Range("P52").Value2 = 208
$ws.Range("Q52").Value2 = 60
